$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Officer Registration Management: register officers one by one onto the
# Pasir Ris Park project (row 5), building up the comma-separated Officer
# list incrementally as each registration is approved.
$ws.Range("N5").Value = "T2109876H"
$ws.Range("N5").Value = "T2109876H, T2109876H"
$ws.Range("N5").Value = "T2109876H, T2109876H, T2109876H"
$ws.Range("N5").Value = "T2109876H, T2109876H, T2109876H, T2109876H"
